$d = $word.ActiveDocument

# 1. Insert a new Heading1 paragraph with the repo URL, followed by an
#    empty paragraph, at the very start of the document body.
$introXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>https://github.com/Saruul-Ulzii/saruul-ulzii.github.io/tree/master/LAB11</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$start = $d.Range(0, 0)
$start.InsertXML($introXml)

# 2. Drop the stale <w:lastRenderedPageBreak/> markers that used to sit in
#    front of the MULTIPLY / SUBTRACT / ADD headings. A self-replace via
#    Find forces Word to re-emit those runs without the cached marker.
$d.Content.Find.Execute("MULTIPLY", $true, $false, $false, $false, $false, $true, 1, $false, "MULTIPLY", 2)
$d.Content.Find.Execute("SUBTRACT", $true, $false, $false, $false, $false, $true, 1, $false, "SUBTRACT", 2)
$d.Content.Find.Execute("ADD", $true, $false, $false, $false, $false, $true, 1, $false, "ADD", 2)
